$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 1126.4
$ws.Range("J2").Value2 = 1642.4286
$ws.Range("L2").Value2 = 1642.4286
$ws.Range("N2").Value2 = -1868.4286
$ws.Range("H43").Value2 = 4213.6665
$ws.Range("J43").Value2 = 3662.1667
$ws.Range("L43").Value2 = 3662.1667
$ws.Range("N43").Value2 = -3800.1667
$ws.Range("H76").Value2 = 12379.4
$ws.Range("I76").Value2 = 25998.5
$ws.Range("K76").Value2 = 25998.5
$ws.Range("M76").Value2 = -25683.5
$ws.Range("H79").Value2 = 12379.4
$ws.Range("I79").Value2 = 25998.5
$ws.Range("K79").Value2 = 25998.5
$ws.Range("M79").Value2 = -24906.5
$ws.Range("H118").Value2 = 4000
$ws.Range("J118").Value2 = 4000
$ws.Range("L118").Value2 = 12000
$ws.Range("N118").Value2 = -15314
$ws.Range("H138").Value2 = 4090.7307
$ws.Range("I138").Value2 = 1647.5652
$ws.Range("J138").Value2 = 22821.666
$ws.Range("K138").Value2 = 4942.6956
$ws.Range("L138").Value2 = 68464.99800000001
$ws.Range("M138").Value2 = 197.3044
$ws.Range("N138").Value2 = -78744.99800000001
$ws.Range("H141").Value2 = 23819074
$ws.Range("I141").Value2 = 31254140
$ws.Range("K141").Value2 = 93762420
$ws.Range("M141").Value2 = -93757240

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 2496.4119
$ws.Range("I2").Value2 = 2587.6924
$ws.Range("K2").Value2 = 2587.6924
$ws.Range("M2").Value2 = -2474.6924
$ws.Range("H32").Value2 = 3877.5833
$ws.Range("I32").Value2 = 3910.3333
$ws.Range("K32").Value2 = 3910.3333
$ws.Range("M32").Value2 = -3623.3333
$ws.Range("H45").Value2 = 1348.5714
$ws.Range("I45").Value2 = 1334.6364
$ws.Range("K45").Value2 = 1334.6364
$ws.Range("M45").Value2 = -957.6364000000001
$ws.Range("H61").Value2 = 17502932
$ws.Range("I61").Value2 = 28574612
$ws.Range("J61").Value2 = 2002580
$ws.Range("K61").Value2 = 28574612
$ws.Range("L61").Value2 = 2002580
$ws.Range("M61").Value2 = -28574400
$ws.Range("N61").Value2 = -2003004
$ws.Range("H74").Value2 = 927660.25
$ws.Range("I74").Value2 = 1042742.2
$ws.Range("K74").Value2 = 1042742.2
$ws.Range("M74").Value2 = -1041868.2
$ws.Range("H77").Value2 = 927660.25
$ws.Range("I77").Value2 = 1042742.2
$ws.Range("K77").Value2 = 5213711
$ws.Range("M77").Value2 = -5209343
$ws.Range("H116").Value2 = 2496.4119
$ws.Range("I116").Value2 = 2587.6924
$ws.Range("K116").Value2 = 2587.6924
$ws.Range("M116").Value2 = -293.6923999999999
$ws.Range("H122").Value2 = 3573.5293
$ws.Range("I122").Value2 = 3483.6667
$ws.Range("J122").Value2 = 3789.2
$ws.Range("K122").Value2 = 10451.0001
$ws.Range("L122").Value2 = 11367.6
$ws.Range("M122").Value2 = -8001.000100000001
$ws.Range("N122").Value2 = -16267.6
$ws.Range("H132").Value2 = 2003587.2
$ws.Range("I132").Value2 = 3612.4211
$ws.Range("J132").Value2 = 8336841
$ws.Range("K132").Value2 = 10837.2633
$ws.Range("L132").Value2 = 25010523
$ws.Range("M132").Value2 = -8307.263300000001
$ws.Range("N132").Value2 = -25015583
$ws.Range("H136").Value2 = 17502932
$ws.Range("I136").Value2 = 28574612
$ws.Range("J136").Value2 = 2002580
$ws.Range("K136").Value2 = 85723836
$ws.Range("L136").Value2 = 6007740
$ws.Range("M136").Value2 = -85721286
$ws.Range("N136").Value2 = -6012840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 2496.4119
$ws.Range("I3").Value2 = 2587.6924
$ws.Range("K3").Value2 = 2587.6924
$ws.Range("M3").Value2 = -2473.6924
$ws.Range("H86").Value2 = 3550.9565
$ws.Range("J86").Value2 = 5273.2
$ws.Range("L86").Value2 = 5273.2
$ws.Range("N86").Value2 = -7519.2
$ws.Range("H89").Value2 = 3550.9565
$ws.Range("J89").Value2 = 5273.2
$ws.Range("L89").Value2 = 26366
$ws.Range("N89").Value2 = -37598
$ws.Range("H99").Value2 = 2547.1667
$ws.Range("I99").Value2 = 2146.7
$ws.Range("K99").Value2 = 2146.7
$ws.Range("M99").Value2 = -648.6999999999998
$ws.Range("H134").Value2 = 4764478.5
$ws.Range("I134").Value2 = 2484.5
$ws.Range("K134").Value2 = 7453.5
$ws.Range("M134").Value2 = -4918.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 43925452
$ws.Range("I31").Value2 = 58827580
$ws.Range("J31").Value2 = 1702756.4
$ws.Range("K31").Value2 = 58827580
$ws.Range("L31").Value2 = 1702756.4
$ws.Range("M31").Value2 = -58827285
$ws.Range("N31").Value2 = -1703346.4
$ws.Range("H34").Value2 = 43925452
$ws.Range("I34").Value2 = 58827580
$ws.Range("J34").Value2 = 1702756.4
$ws.Range("K34").Value2 = 58827580
$ws.Range("L34").Value2 = 1702756.4
$ws.Range("M34").Value2 = -58827378
$ws.Range("N34").Value2 = -1703160.4
$ws.Range("H105").Value2 = 3224.5
$ws.Range("I105").Value2 = 1450
$ws.Range("K105").Value2 = 1450
$ws.Range("M105").Value2 = 297
$ws.Range("H107").Value2 = 3430.879
$ws.Range("I107").Value2 = 3070.7827
$ws.Range("K107").Value2 = 3070.7827
$ws.Range("M107").Value2 = -1150.7827
$ws.Range("H134").Value2 = 2152.2222
$ws.Range("I134").Value2 = 2213.1428
$ws.Range("J134").Value2 = 1939
$ws.Range("K134").Value2 = 6639.428400000001
$ws.Range("L134").Value2 = 5817
$ws.Range("M134").Value2 = -4104.428400000001
$ws.Range("N134").Value2 = -10887

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 151.27272
$ws.Range("J2").Value2 = 177.55556
$ws.Range("L2").Value2 = 1065.33336
$ws.Range("N2").Value2 = -1291.33336
$ws.Range("H124").Value2 = 38828.832
$ws.Range("I124").Value2 = 44353.332
$ws.Range("K124").Value2 = 133059.996
$ws.Range("M124").Value2 = -128149.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 842.9167
$ws.Range("I97").Value2 = 678.8
$ws.Range("J97").Value2 = 960.1429000000001
$ws.Range("K97").Value2 = 678.8
$ws.Range("L97").Value2 = 960.1429000000001
$ws.Range("M97").Value2 = -182.8
$ws.Range("N97").Value2 = -1952.1429
$ws.Range("H102").Value2 = 2733.1667
$ws.Range("I102").Value2 = 2718.182
$ws.Range("K102").Value2 = 2718.182
$ws.Range("M102").Value2 = -1096.182
$ws.Range("H113").Value2 = 928494.25
$ws.Range("J113").Value2 = 3088585.2
$ws.Range("L113").Value2 = 3088585.2
$ws.Range("N113").Value2 = -3092925.2
$ws.Range("H122").Value2 = 4916.273
$ws.Range("I122").Value2 = 6872
$ws.Range("K122").Value2 = 20616
$ws.Range("M122").Value2 = -18166
$ws.Range("H132").Value2 = 3574397.5
$ws.Range("I132").Value2 = 2868.7144
$ws.Range("J132").Value2 = 14288984
$ws.Range("K132").Value2 = 8606.143199999999
$ws.Range("L132").Value2 = 42866952
$ws.Range("M132").Value2 = -6076.143199999999
$ws.Range("N132").Value2 = -42872012

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 3589
$ws.Range("I16").Value2 = 2034.5625
$ws.Range("K16").Value2 = 2034.5625
$ws.Range("M16").Value2 = -1864.5625
$ws.Range("H22").Value2 = 8857.210999999999
$ws.Range("I22").Value2 = 11022.462
$ws.Range("J22").Value2 = 4165.8335
$ws.Range("K22").Value2 = 11022.462
$ws.Range("L22").Value2 = 4165.8335
$ws.Range("M22").Value2 = -10727.462
$ws.Range("N22").Value2 = -4755.8335
$ws.Range("H27").Value2 = 8857.210999999999
$ws.Range("I27").Value2 = 11022.462
$ws.Range("J27").Value2 = 4165.8335
$ws.Range("K27").Value2 = 11022.462
$ws.Range("L27").Value2 = 4165.8335
$ws.Range("M27").Value2 = -10915.462
$ws.Range("N27").Value2 = -4379.8335
$ws.Range("H46").Value2 = 6699.6
$ws.Range("I46").Value2 = 6749
$ws.Range("J46").Value2 = 6666.6665
$ws.Range("K46").Value2 = 6749
$ws.Range("L46").Value2 = 6666.6665
$ws.Range("M46").Value2 = -6561
$ws.Range("N46").Value2 = -7042.6665
$ws.Range("H97").Value2 = 52499.832
$ws.Range("J97").Value2 = 52499.832
$ws.Range("L97").Value2 = 52499.832
$ws.Range("N97").Value2 = -54481.832
$ws.Range("H136").Value2 = 3619.9473
$ws.Range("J136").Value2 = 4110.8887
$ws.Range("L136").Value2 = 12332.6661
$ws.Range("N136").Value2 = -17432.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value2 = 10750
$ws.Range("J4").Value2 = 10000
$ws.Range("L4").Value2 = 10000
$ws.Range("N4").Value2 = -10226
$ws.Range("H136").Value2 = 279413.38
$ws.Range("I136").Value2 = 1615.9697
$ws.Range("K136").Value2 = 4847.909100000001
$ws.Range("M136").Value2 = -2297.909100000001

Write-Output "Applied all Ragnarok_Profits updates"